$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark from the first paragraph
#    (it currently sits right after "... DESIGNING FOR EMERGING TECHNOLOGIES").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Delete the four trailing body paragraphs entirely
#    ("McElroy discusses...", "In prototypes for...", "The material aspect...",
#     "When prototyping for...") -- these are paragraphs 4 through 7.
$deleteRange = $d.Range($d.Paragraphs.Item(4).Range.Start, $d.Paragraphs.Item(7).Range.End)
$deleteRange.Delete()

# 3. Replace the whole content of paragraph 3 (the two runs making up the
#    "In her text, McElroy begins..." / "The category I find..." paragraph)
#    with a temporary placeholder "AXX" so we can safely anchor a bookmark
#    one character in (avoiding an edge case at the very end of the
#    paragraph), then trim back down to just "A".
$p3 = $d.Paragraphs.Item(3)
$bodyRange = $p3.Range
$bodyRange.MoveEnd(1, -1)
$bodyRange.Text = "AXX"

$p3Start = $d.Paragraphs.Item(3).Range.Start
$bookmarkPos = $p3Start + 1

# 4. Re-create the "_GoBack" bookmark collapsed right after the "A", then
#    trim off the temporary "XX" placeholder text that followed it.
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$trimRange = $d.Range($bookmarkPos, $bookmarkPos + 2)
$trimRange.Delete()
